$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 374).
$ws.Range("C2:C374").Value = 46061
